$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 4848.2856
$ws.Cells.Item(76, 9).Value = 3815.6667
$ws.Cells.Item(76, 10).Value = 5622.75
$ws.Cells.Item(76, 11).Value = 3815.6667
$ws.Cells.Item(76, 12).Value = 5622.75
$ws.Cells.Item(76, 13).Value = -3500.6667
$ws.Cells.Item(76, 14).Value = -6252.75

$ws.Cells.Item(79, 8).Value = 4848.2856
$ws.Cells.Item(79, 9).Value = 3815.6667
$ws.Cells.Item(79, 10).Value = 5622.75
$ws.Cells.Item(79, 11).Value = 3815.6667
$ws.Cells.Item(79, 12).Value = 5622.75
$ws.Cells.Item(79, 13).Value = -2723.6667
$ws.Cells.Item(79, 14).Value = -7806.75

$ws.Cells.Item(138, 8).Value = 3501.8625
$ws.Cells.Item(138, 9).Value = 1830.5264
$ws.Cells.Item(138, 10).Value = 4022.4426
$ws.Cells.Item(138, 11).Value = 5491.5792
$ws.Cells.Item(138, 12).Value = 12067.3278
$ws.Cells.Item(138, 13).Value = -351.5792000000001
$ws.Cells.Item(138, 14).Value = -22347.3278

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11322.667
$ws.Cells.Item(32, 9).Value = 10752.404
$ws.Cells.Item(32, 10).Value = 15314.5
$ws.Cells.Item(32, 11).Value = 10752.404
$ws.Cells.Item(32, 12).Value = 15314.5
$ws.Cells.Item(32, 13).Value = -10465.404
$ws.Cells.Item(32, 14).Value = -15888.5

$ws.Cells.Item(61, 8).Value = 7423.4316
$ws.Cells.Item(61, 9).Value = 7849.4688
$ws.Cells.Item(61, 11).Value = 7849.4688
$ws.Cells.Item(61, 13).Value = -7637.4688

$ws.Cells.Item(136, 8).Value = 7423.4316
$ws.Cells.Item(136, 9).Value = 7849.4688
$ws.Cells.Item(136, 11).Value = 23548.4064
$ws.Cells.Item(136, 13).Value = -20998.4064

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2420.7058
$ws.Cells.Item(94, 9).Value = 1255.16
$ws.Cells.Item(94, 11).Value = 1255.16
$ws.Cells.Item(94, 13).Value = -804.1600000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 223.10527
$ws.Cells.Item(7, 9).Value = 216.72728
$ws.Cells.Item(7, 10).Value = 231.875
$ws.Cells.Item(7, 11).Value = 216.72728
$ws.Cells.Item(7, 12).Value = 231.875
$ws.Cells.Item(7, 13).Value = -103.72728
$ws.Cells.Item(7, 14).Value = -457.875

$ws.Cells.Item(86, 8).Value = 6156.7617
$ws.Cells.Item(86, 9).Value = 4965.125
$ws.Cells.Item(86, 10).Value = 9970
$ws.Cells.Item(86, 11).Value = 4965.125
$ws.Cells.Item(86, 12).Value = 9970
$ws.Cells.Item(86, 13).Value = -3842.125
$ws.Cells.Item(86, 14).Value = -12216

$ws.Cells.Item(89, 8).Value = 6156.7617
$ws.Cells.Item(89, 9).Value = 4965.125
$ws.Cells.Item(89, 10).Value = 9970
$ws.Cells.Item(89, 11).Value = 24825.625
$ws.Cells.Item(89, 12).Value = 49850
$ws.Cells.Item(89, 13).Value = -19209.625
$ws.Cells.Item(89, 14).Value = -61082

$ws.Cells.Item(99, 8).Value = 3230298.8
$ws.Cells.Item(99, 9).Value = 5282888.5
$ws.Cells.Item(99, 11).Value = 5282888.5
$ws.Cells.Item(99, 13).Value = -5281390.5

$ws.Cells.Item(122, 8).Value = 11938
$ws.Cells.Item(122, 9).Value = 18780.572
$ws.Cells.Item(122, 11).Value = 56341.716
$ws.Cells.Item(122, 13).Value = -53891.716

$ws.Cells.Item(126, 8).Value = 3230298.8
$ws.Cells.Item(126, 9).Value = 5282888.5
$ws.Cells.Item(126, 11).Value = 15848665.5
$ws.Cells.Item(126, 13).Value = -15846195.5

$ws.Cells.Item(132, 8).Value = 1932
$ws.Cells.Item(132, 9).Value = 1727.5
$ws.Cells.Item(132, 11).Value = 5182.5
$ws.Cells.Item(132, 13).Value = -2652.5

$ws.Cells.Item(141, 8).Value = 178736.23
$ws.Cells.Item(141, 10).Value = 188831.1
$ws.Cells.Item(141, 12).Value = 188831.1
$ws.Cells.Item(141, 14).Value = -199191.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 975.94116
$ws.Cells.Item(39, 10).Value = 12000
$ws.Cells.Item(39, 12).Value = 36000
$ws.Cells.Item(39, 14).Value = -36588

$ws.Cells.Item(121, 8).Value = 793119.5
$ws.Cells.Item(121, 9).Value = 1408279.1
$ws.Cells.Item(121, 10).Value = 2200
$ws.Cells.Item(121, 11).Value = 4224837.300000001
$ws.Cells.Item(121, 12).Value = 6600
$ws.Cells.Item(121, 13).Value = -4223527.300000001
$ws.Cells.Item(121, 14).Value = -9220

$ws.Cells.Item(134, 8).Value = 13803.258
$ws.Cells.Item(134, 9).Value = 24136.072
$ws.Cells.Item(134, 10).Value = 5293.8823
$ws.Cells.Item(134, 11).Value = 72408.216
$ws.Cells.Item(134, 12).Value = 15881.6469
$ws.Cells.Item(134, 13).Value = -67338.216
$ws.Cells.Item(134, 14).Value = -26021.6469

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 18069.727
$ws.Cells.Item(80, 9).Value = 19888
$ws.Cells.Item(80, 11).Value = 19888
$ws.Cells.Item(80, 13).Value = -18890

$ws.Cells.Item(83, 8).Value = 18069.727
$ws.Cells.Item(83, 9).Value = 19888
$ws.Cells.Item(83, 11).Value = 99440
$ws.Cells.Item(83, 13).Value = -94448

$ws.Cells.Item(122, 8).Value = 9083.869000000001
$ws.Cells.Item(122, 9).Value = 6808.125
$ws.Cells.Item(122, 11).Value = 20424.375
$ws.Cells.Item(122, 13).Value = -17974.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 16289.6
$ws.Cells.Item(7, 9).Value = 21996.227
$ws.Cells.Item(7, 11).Value = 21996.227
$ws.Cells.Item(7, 13).Value = -21884.227

$ws.Cells.Item(38, 8).Value = 49998
$ws.Cells.Item(38, 10).Value = 49998
$ws.Cells.Item(38, 12).Value = 49998
$ws.Cells.Item(38, 14).Value = -50818

$ws.Cells.Item(68, 8).Value = 4899.1763
$ws.Cells.Item(68, 9).Value = 2040.1666
$ws.Cells.Item(68, 10).Value = 6458.636
$ws.Cells.Item(68, 11).Value = 2040.1666
$ws.Cells.Item(68, 12).Value = 6458.636
$ws.Cells.Item(68, 13).Value = -1291.1666
$ws.Cells.Item(68, 14).Value = -7956.636

$ws.Cells.Item(71, 8).Value = 4899.1763
$ws.Cells.Item(71, 9).Value = 2040.1666
$ws.Cells.Item(71, 10).Value = 6458.636
$ws.Cells.Item(71, 11).Value = 10200.833
$ws.Cells.Item(71, 12).Value = 32293.18
$ws.Cells.Item(71, 13).Value = -6456.833000000001
$ws.Cells.Item(71, 14).Value = -39781.18

$ws.Cells.Item(126, 8).Value = 16289.6
$ws.Cells.Item(126, 9).Value = 21996.227
$ws.Cells.Item(126, 11).Value = 65988.681
$ws.Cells.Item(126, 13).Value = -63518.681

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 6332.579
$ws.Cells.Item(81, 9).Value = 9469.833000000001
$ws.Cells.Item(81, 10).Value = 954.4286
$ws.Cells.Item(81, 11).Value = 18939.666
$ws.Cells.Item(81, 12).Value = 1908.8572
$ws.Cells.Item(81, 13).Value = -17878.666
$ws.Cells.Item(81, 14).Value = -4030.8572

$ws.Cells.Item(84, 8).Value = 6332.579
$ws.Cells.Item(84, 9).Value = 9469.833000000001
$ws.Cells.Item(84, 10).Value = 954.4286
$ws.Cells.Item(84, 11).Value = 94698.33
$ws.Cells.Item(84, 12).Value = 9544.286
$ws.Cells.Item(84, 13).Value = -89394.33
$ws.Cells.Item(84, 14).Value = -20152.286

$ws.Cells.Item(107, 8).Value = 40237.125
$ws.Cells.Item(107, 9).Value = 3599.4
$ws.Cells.Item(107, 11).Value = 10798.2
$ws.Cells.Item(107, 13).Value = -8878.200000000001

$ws.Cells.Item(136, 8).Value = 291935.3
$ws.Cells.Item(136, 9).Value = 297434.06
$ws.Cells.Item(136, 10).Value = 5999
$ws.Cells.Item(136, 11).Value = 892302.1799999999
$ws.Cells.Item(136, 12).Value = 17997
$ws.Cells.Item(136, 13).Value = -889752.1799999999
$ws.Cells.Item(136, 14).Value = -23097
